$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Arkusz1")
$ws.Range("H284").Value = "test"
